# daily auto push: 2026-01-18 13:37 UTC
#
# A new reading was logged for 2026/01/18 (day-of-week "日") at slot 19,
# ranking 201. It belongs right after the existing 2026/01/18 rows
# (row 664, value 16) and before the 2026/12/29 block that currently
# starts at row 664 - so insert a new row there and push everything
# else down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 664..end down by one to make room for the new entry.
$ws.Rows.Item(664).Insert()

# Column A holds dates stored as literal text (e.g. "2026/01/18"), not
# real date serials. Setting .Value directly on a General-formatted
# cell makes Excel's COM layer auto-parse the date-like string into a
# date serial, which we don't want. Mark the cell as Text first so the
# string is stored verbatim, then copy the (unstyled) format from the
# row above back onto it so we don't leave a stray number-format behind.
$ws.Range("A664").NumberFormat = "@"
$ws.Range("A664").Value = "2026/01/18"
$ws.Range("A663").Copy()
$ws.Range("A664").PasteSpecial(-4122)

$ws.Range("B664").Value = "日"
$ws.Range("C664").Value = 19
$ws.Range("D664").Value = 201
